$wb = $excel.ActiveWorkbook

# Rename sheets: rename sheet 2 first to avoid name collisions,
# then rename sheet 1. sheetId/rId order is preserved (no sheet add/remove).
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Marzo"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Febrero"

# Update data values for sheet that is now "Febrero" (was "Enero", rId1)
$ws1.Range("B2").Value = 20207
$ws1.Range("C2").Value = 748
$ws1.Range("D2").Value = 470
$ws1.Range("E2").Value = 1377
$ws1.Range("F2").Value = 17612
$ws1.Range("G2").Value = 16.28
$ws1.Range("H2").Value = 18.05
$ws1.Range("I2").Value = 0
$ws1.Range("J2").Value = 0
$ws1.Range("K2").Value = 1
$ws1.Range("B3").Value = 14599
$ws1.Range("C3").Value = 2857
$ws1.Range("D3").Value = 636
$ws1.Range("E3").Value = 3186
$ws1.Range("F3").Value = 7920
$ws1.Range("G3").Value = 7.69
$ws1.Range("H3").Value = 5.74
$ws1.Range("I3").Value = 0
$ws1.Range("J3").Value = 5
$ws1.Range("K3").Value = 8
$ws1.Range("B4").Value = 7823
$ws1.Range("C4").Value = 542
$ws1.Range("D4").Value = 321
$ws1.Range("E4").Value = 976
$ws1.Range("F4").Value = 5984
$ws1.Range("G4").Value = 8.359999999999999
$ws1.Range("H4").Value = 6.77
$ws1.Range("I4").Value = 4
$ws1.Range("J4").Value = 2
$ws1.Range("K4").Value = 6
$ws1.Range("B5").Value = 25013
$ws1.Range("C5").Value = 187
$ws1.Range("D5").Value = 1342
$ws1.Range("E5").Value = 2814
$ws1.Range("F5").Value = 20670
$ws1.Range("G5").Value = 20.64
$ws1.Range("H5").Value = 21.9
$ws1.Range("I5").Value = 0
$ws1.Range("J5").Value = 0
$ws1.Range("K5").Value = 2
$ws1.Range("B6").Value = 13541
$ws1.Range("C6").Value = 322
$ws1.Range("D6").Value = 523
$ws1.Range("E6").Value = 1751
$ws1.Range("F6").Value = 10945
$ws1.Range("G6").Value = 3.97
$ws1.Range("H6").Value = 5.35
$ws1.Range("I6").Value = 1
$ws1.Range("J6").Value = 8
$ws1.Range("K6").Value = 5
$ws1.Range("B7").Value = 10824
$ws1.Range("C7").Value = 1947
$ws1.Range("D7").Value = 681
$ws1.Range("E7").Value = 1146
$ws1.Range("F7").Value = 7050
$ws1.Range("G7").Value = 24
$ws1.Range("H7").Value = 9.25
$ws1.Range("I7").Value = 18
$ws1.Range("J7").Value = 13
$ws1.Range("K7").Value = 12
$ws1.Range("B8").Value = 13533
$ws1.Range("C8").Value = 1459
$ws1.Range("D8").Value = 878
$ws1.Range("E8").Value = 1203
$ws1.Range("F8").Value = 9993
$ws1.Range("G8").Value = 10.2
$ws1.Range("H8").Value = 13.66
$ws1.Range("I8").Value = 5
$ws1.Range("J8").Value = 9
$ws1.Range("K8").Value = 5
$ws1.Range("B9").Value = 16415
$ws1.Range("C9").Value = 136
$ws1.Range("D9").Value = 280
$ws1.Range("E9").Value = 566
$ws1.Range("F9").Value = 15433
$ws1.Range("G9").Value = 18.86
$ws1.Range("H9").Value = 5.6
$ws1.Range("I9").Value = 0
$ws1.Range("J9").Value = 0
$ws1.Range("K9").Value = 4
$ws1.Range("B10").Value = 2645
$ws1.Range("C10").Value = 886
$ws1.Range("D10").Value = 370
$ws1.Range("E10").Value = 1389
$ws1.Range("F10").Value = 0
$ws1.Range("G10").Value = 4.69
$ws1.Range("H10").Value = 3.36
$ws1.Range("I10").Value = 13
$ws1.Range("J10").Value = 11
$ws1.Range("K10").Value = 8
$ws1.Range("B11").Value = 15789
$ws1.Range("C11").Value = 563
$ws1.Range("D11").Value = 247
$ws1.Range("E11").Value = 696
$ws1.Range("F11").Value = 14283
$ws1.Range("G11").Value = 8.43
$ws1.Range("H11").Value = 8.140000000000001
$ws1.Range("I11").Value = 5
$ws1.Range("J11").Value = 9
$ws1.Range("K11").Value = 0
$ws1.Range("B12").Value = 8706
$ws1.Range("C12").Value = 1208
$ws1.Range("D12").Value = 562
$ws1.Range("E12").Value = 1178
$ws1.Range("F12").Value = 5756
$ws1.Range("G12").Value = 11.81
$ws1.Range("H12").Value = 8.26
$ws1.Range("I12").Value = 40
$ws1.Range("J12").Value = 35
$ws1.Range("K12").Value = 31
$ws1.Range("B13").Value = 18114
$ws1.Range("C13").Value = 357
$ws1.Range("D13").Value = 704
$ws1.Range("E13").Value = 1753
$ws1.Range("F13").Value = 15299
$ws1.Range("G13").Value = 11.01
$ws1.Range("H13").Value = 11.8
$ws1.Range("I13").Value = 6
$ws1.Range("J13").Value = 17
$ws1.Range("K13").Value = 7
$ws1.Range("B14").Value = 18311
$ws1.Range("C14").Value = 442
$ws1.Range("D14").Value = 375
$ws1.Range("E14").Value = 971
$ws1.Range("F14").Value = 16522
$ws1.Range("G14").Value = 17.57
$ws1.Range("H14").Value = 11.83
$ws1.Range("I14").Value = 0
$ws1.Range("J14").Value = 0
$ws1.Range("K14").Value = 5
$ws1.Range("B15").Value = 14599
$ws1.Range("C15").Value = 2857
$ws1.Range("D15").Value = 636
$ws1.Range("E15").Value = 3186
$ws1.Range("F15").Value = 7920
$ws1.Range("G15").Value = 7.69
$ws1.Range("H15").Value = 5.74
$ws1.Range("I15").Value = 0
$ws1.Range("J15").Value = 5
$ws1.Range("K15").Value = 8

# Update data values for sheet that is now "Marzo" (was "Febrero", rId2)
$ws2.Range("B2").Value = 15356
$ws2.Range("C2").Value = 2398
$ws2.Range("D2").Value = 669
$ws2.Range("E2").Value = 894
$ws2.Range("F2").Value = 11395
$ws2.Range("G2").Value = 27.85
$ws2.Range("H2").Value = 9.67
$ws2.Range("I2").Value = 0
$ws2.Range("J2").Value = 0
$ws2.Range("K2").Value = 9
$ws2.Range("B3").Value = 17551
$ws2.Range("C3").Value = 3395
$ws2.Range("D3").Value = 1973
$ws2.Range("E3").Value = 1442
$ws2.Range("F3").Value = 10741
$ws2.Range("G3").Value = 3.52
$ws2.Range("H3").Value = 24.93
$ws2.Range("I3").Value = 0
$ws2.Range("J3").Value = 8
$ws2.Range("K3").Value = 6
$ws2.Range("B4").Value = 13204
$ws2.Range("C4").Value = 870
$ws2.Range("D4").Value = 193
$ws2.Range("E4").Value = 369
$ws2.Range("F4").Value = 11772
$ws2.Range("G4").Value = 6.42
$ws2.Range("H4").Value = 4.28
$ws2.Range("I4").Value = 5
$ws2.Range("J4").Value = 2
$ws2.Range("K4").Value = 11
$ws2.Range("B5").Value = 19823
$ws2.Range("C5").Value = 1820
$ws2.Range("D5").Value = 917
$ws2.Range("E5").Value = 2963
$ws2.Range("F5").Value = 14123
$ws2.Range("G5").Value = 3.95
$ws2.Range("H5").Value = 9.57
$ws2.Range("I5").Value = 2
$ws2.Range("J5").Value = 0
$ws2.Range("K5").Value = 2
$ws2.Range("B6").Value = 11574
$ws2.Range("C6").Value = 418
$ws2.Range("D6").Value = 614
$ws2.Range("E6").Value = 3499
$ws2.Range("F6").Value = 7043
$ws2.Range("G6").Value = 9.15
$ws2.Range("H6").Value = 3.92
$ws2.Range("I6").Value = 3
$ws2.Range("J6").Value = 6
$ws2.Range("K6").Value = 9
$ws2.Range("B7").Value = 9131
$ws2.Range("C7").Value = 3026
$ws2.Range("D7").Value = 672
$ws2.Range("E7").Value = 1370
$ws2.Range("F7").Value = 4063
$ws2.Range("G7").Value = 3.18
$ws2.Range("H7").Value = 1.84
$ws2.Range("I7").Value = 19
$ws2.Range("J7").Value = 18
$ws2.Range("K7").Value = 2
$ws2.Range("B8").Value = 19179
$ws2.Range("C8").Value = 467
$ws2.Range("D8").Value = 1579
$ws2.Range("E8").Value = 2255
$ws2.Range("F8").Value = 14878
$ws2.Range("G8").Value = 16.52
$ws2.Range("H8").Value = 6.99
$ws2.Range("I8").Value = 7
$ws2.Range("J8").Value = 9
$ws2.Range("K8").Value = 6
$ws2.Range("B9").Value = 19363
$ws2.Range("C9").Value = 718
$ws2.Range("D9").Value = 725
$ws2.Range("E9").Value = 259
$ws2.Range("F9").Value = 17661
$ws2.Range("G9").Value = 9.029999999999999
$ws2.Range("H9").Value = 5.67
$ws2.Range("I9").Value = 0
$ws2.Range("J9").Value = 0
$ws2.Range("K9").Value = 1
$ws2.Range("B10").Value = 19287
$ws2.Range("C10").Value = 1064
$ws2.Range("D10").Value = 856
$ws2.Range("E10").Value = 596
$ws2.Range("F10").Value = 16771
$ws2.Range("G10").Value = 7.15
$ws2.Range("H10").Value = 12.95
$ws2.Range("I10").Value = 21
$ws2.Range("J10").Value = 12
$ws2.Range("K10").Value = 6
$ws2.Range("B11").Value = 6203
$ws2.Range("C11").Value = 2297
$ws2.Range("D11").Value = 446
$ws2.Range("E11").Value = 564
$ws2.Range("F11").Value = 2896
$ws2.Range("G11").Value = 10.82
$ws2.Range("H11").Value = 9.380000000000001
$ws2.Range("I11").Value = 6
$ws2.Range("J11").Value = 21
$ws2.Range("K11").Value = 0
$ws2.Range("B12").Value = 15200
$ws2.Range("C12").Value = 1356
$ws2.Range("D12").Value = 825
$ws2.Range("E12").Value = 1147
$ws2.Range("F12").Value = 11871
$ws2.Range("G12").Value = 8.32
$ws2.Range("H12").Value = 6.52
$ws2.Range("I12").Value = 52
$ws2.Range("J12").Value = 41
$ws2.Range("K12").Value = 25
$ws2.Range("B13").Value = 12533
$ws2.Range("C13").Value = 1511
$ws2.Range("D13").Value = 659
$ws2.Range("E13").Value = 2342
$ws2.Range("F13").Value = 8020
$ws2.Range("G13").Value = 7.97
$ws2.Range("H13").Value = 7.62
$ws2.Range("I13").Value = 11
$ws2.Range("J13").Value = 27
$ws2.Range("K13").Value = 11
$ws2.Range("B14").Value = 17359
$ws2.Range("C14").Value = 1558
$ws2.Range("D14").Value = 697
$ws2.Range("E14").Value = 576
$ws2.Range("F14").Value = 14528
$ws2.Range("G14").Value = 18.44
$ws2.Range("H14").Value = 7.66
$ws2.Range("I14").Value = 0
$ws2.Range("J14").Value = 0
$ws2.Range("K14").Value = 10
$ws2.Range("B15").Value = 17551
$ws2.Range("C15").Value = 3395
$ws2.Range("D15").Value = 1973
$ws2.Range("E15").Value = 1442
$ws2.Range("F15").Value = 10741
$ws2.Range("G15").Value = 3.52
$ws2.Range("H15").Value = 24.93
$ws2.Range("I15").Value = 0
$ws2.Range("J15").Value = 8
$ws2.Range("K15").Value = 6

Write-Output "done"
